# Remove the "Phụ cấp" (allowance) rows for CẦN THƠ and LONG XUYÊN from the
# "Lương" sheet. The "Phụ cấp tại SÓC TRĂNG" row is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

# Row 3 is "Phụ cấp tại CẦN THƠ" - delete it; everything below shifts up.
$ws.Rows.Item(3).Delete()

# After the row-3 deletion, "Phụ cấp tại LONG XUYÊN" (originally row 14)
# is now at row 13 - delete it too.
$ws.Rows.Item(13).Delete()
